$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table 1 (rows 3-12): Standart AWM results
$ws.Cells.Item(3, 1).Value = 0
$ws.Cells.Item(3, 2).Value = 87
$ws.Cells.Item(3, 3).Value = 'Praia do Farol'
$ws.Cells.Item(3, 4).Value = 'beach'
$ws.Cells.Item(3, 5).Value = 1.581773681253018
$ws.Cells.Item(3, 6).Value = 0.791
$ws.Cells.Item(3, 7).Value = -12.952961
$ws.Cells.Item(3, 8).Value = -38.348126

$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = 193
$ws.Cells.Item(4, 3).Value = 'Muncab - National Museum of Afro Brazilian Culture'
$ws.Cells.Item(4, 4).Value = 'museum'
$ws.Cells.Item(4, 5).Value = 1.442942753776769
$ws.Cells.Item(4, 6).Value = 0.746
$ws.Cells.Item(4, 7).Value = -12.9761485
$ws.Cells.Item(4, 8).Value = -38.5124649

$ws.Cells.Item(5, 1).Value = 2
$ws.Cells.Item(5, 2).Value = 77
$ws.Cells.Item(5, 3).Value = 'Praia do Porto da Barra'
$ws.Cells.Item(5, 4).Value = 'beach'
$ws.Cells.Item(5, 5).Value = 1.484756941195069
$ws.Cells.Item(5, 6).Value = 0.742
$ws.Cells.Item(5, 7).Value = -13.0038021
$ws.Cells.Item(5, 8).Value = -38.5326932

$ws.Cells.Item(6, 1).Value = 3
$ws.Cells.Item(6, 2).Value = 199
$ws.Cells.Item(6, 3).Value = 'Museum of Brazilian Music'
$ws.Cells.Item(6, 4).Value = 'museum'
$ws.Cells.Item(6, 5).Value = 1.409365989240249
$ws.Cells.Item(6, 6).Value = 0.729
$ws.Cells.Item(6, 7).Value = -12.9738253
$ws.Cells.Item(6, 8).Value = -38.508487

$ws.Cells.Item(7, 1).Value = 4
$ws.Cells.Item(7, 2).Value = 349
$ws.Cells.Item(7, 3).Value = 'Shopping Capemi Salvador'
$ws.Cells.Item(7, 4).Value = 'shopping_mall'
$ws.Cells.Item(7, 5).Value = 1.353198932593896
$ws.Cells.Item(7, 6).Value = 0.714
$ws.Cells.Item(7, 7).Value = -12.9831
$ws.Cells.Item(7, 8).Value = -38.465

$ws.Cells.Item(8, 1).Value = 5
$ws.Cells.Item(8, 2).Value = 352
$ws.Cells.Item(8, 3).Value = 'Cia Maritima Salvador Shopping'
$ws.Cells.Item(8, 4).Value = 'shopping_mall'
$ws.Cells.Item(8, 5).Value = 1.279649308173254
$ws.Cells.Item(8, 6).Value = 0.677
$ws.Cells.Item(8, 7).Value = -12.9954162
$ws.Cells.Item(8, 8).Value = -38.4542581

$ws.Cells.Item(9, 1).Value = 6
$ws.Cells.Item(9, 2).Value = 378
$ws.Cells.Item(9, 3).Value = 'Forte de Santo Antonio alem do Carmo'
$ws.Cells.Item(9, 4).Value = 'tourist_attraction'
$ws.Cells.Item(9, 5).Value = 1.33246946639442
$ws.Cells.Item(9, 6).Value = 0.673
$ws.Cells.Item(9, 7).Value = -12.9636017
$ws.Cells.Item(9, 8).Value = -38.5039186

$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = 342
$ws.Cells.Item(10, 3).Value = 'Shopping Paralela'
$ws.Cells.Item(10, 4).Value = 'shopping_mall'
$ws.Cells.Item(10, 5).Value = 1.246563807333256
$ws.Cells.Item(10, 6).Value = 0.661
$ws.Cells.Item(10, 7).Value = -12.9365263
$ws.Cells.Item(10, 8).Value = -38.3949236

$ws.Cells.Item(11, 1).Value = 8
$ws.Cells.Item(11, 2).Value = 83
$ws.Cells.Item(11, 3).Value = 'Ondina Beach'
$ws.Cells.Item(11, 4).Value = 'beach'
$ws.Cells.Item(11, 5).Value = 1.245057115085724
$ws.Cells.Item(11, 6).Value = 0.66
$ws.Cells.Item(11, 7).Value = -13.0103784
$ws.Cells.Item(11, 8).Value = -38.5120368

$ws.Cells.Item(12, 1).Value = 9
$ws.Cells.Item(12, 2).Value = 341
$ws.Cells.Item(12, 3).Value = 'Shopping da Bahia'
$ws.Cells.Item(12, 4).Value = 'shopping_mall'
$ws.Cells.Item(12, 5).Value = 1.246334170797274
$ws.Cells.Item(12, 6).Value = 0.66
$ws.Cells.Item(12, 7).Value = -12.9811659
$ws.Cells.Item(12, 8).Value = -38.46503209999999

$ws.Cells.Item(13, 1).Value = 'NDCG: 1.0'
$ws.Cells.Item(15, 1).Value = 'Diversificado_recs_greedy AWM'

# Table 2 (rows 17-26): Diversificado_recs_greedy AWM results
$ws.Cells.Item(17, 1).Value = 0
$ws.Cells.Item(17, 2).Value = 87
$ws.Cells.Item(17, 3).Value = 'Praia do Farol'
$ws.Cells.Item(17, 4).Value = 'beach'
$ws.Cells.Item(17, 5).Value = 1.581773681253018
$ws.Cells.Item(17, 6).Value = 0.791
$ws.Cells.Item(17, 7).Value = -12.952961
$ws.Cells.Item(17, 8).Value = -38.348126

$ws.Cells.Item(18, 1).Value = 1
$ws.Cells.Item(18, 2).Value = 193
$ws.Cells.Item(18, 3).Value = 'Muncab - National Museum of Afro Brazilian Culture'
$ws.Cells.Item(18, 4).Value = 'museum'
$ws.Cells.Item(18, 5).Value = 1.442942753776769
$ws.Cells.Item(18, 6).Value = 0.746
$ws.Cells.Item(18, 7).Value = -12.9761485
$ws.Cells.Item(18, 8).Value = -38.5124649

$ws.Cells.Item(19, 1).Value = 2
$ws.Cells.Item(19, 2).Value = 349
$ws.Cells.Item(19, 3).Value = 'Shopping Capemi Salvador'
$ws.Cells.Item(19, 4).Value = 'shopping_mall'
$ws.Cells.Item(19, 5).Value = 1.353198932593896
$ws.Cells.Item(19, 6).Value = 0.714
$ws.Cells.Item(19, 7).Value = -12.9831
$ws.Cells.Item(19, 8).Value = -38.465

$ws.Cells.Item(20, 1).Value = 3
$ws.Cells.Item(20, 2).Value = 378
$ws.Cells.Item(20, 3).Value = 'Forte de Santo Antonio alem do Carmo'
$ws.Cells.Item(20, 4).Value = 'tourist_attraction'
$ws.Cells.Item(20, 5).Value = 1.33246946639442
$ws.Cells.Item(20, 6).Value = 0.673
$ws.Cells.Item(20, 7).Value = -12.9636017
$ws.Cells.Item(20, 8).Value = -38.5039186

$ws.Cells.Item(21, 1).Value = 4
$ws.Cells.Item(21, 2).Value = 140
$ws.Cells.Item(21, 3).Value = 'CAFETERIA FLORIDA'
$ws.Cells.Item(21, 4).Value = 'cafe'
$ws.Cells.Item(21, 5).Value = 1.264779587120852
$ws.Cells.Item(21, 6).Value = 0.654
$ws.Cells.Item(21, 7).Value = -13.0039875
$ws.Cells.Item(21, 8).Value = -38.5023199

$ws.Cells.Item(22, 1).Value = 5
$ws.Cells.Item(22, 2).Value = 77
$ws.Cells.Item(22, 3).Value = 'Praia do Porto da Barra'
$ws.Cells.Item(22, 4).Value = 'beach'
$ws.Cells.Item(22, 5).Value = 1.484756941195069
$ws.Cells.Item(22, 6).Value = 0.742
$ws.Cells.Item(22, 7).Value = -13.0038021
$ws.Cells.Item(22, 8).Value = -38.5326932

$ws.Cells.Item(23, 1).Value = 6
$ws.Cells.Item(23, 2).Value = 199
$ws.Cells.Item(23, 3).Value = 'Museum of Brazilian Music'
$ws.Cells.Item(23, 4).Value = 'museum'
$ws.Cells.Item(23, 5).Value = 1.409365989240249
$ws.Cells.Item(23, 6).Value = 0.729
$ws.Cells.Item(23, 7).Value = -12.9738253
$ws.Cells.Item(23, 8).Value = -38.508487

$ws.Cells.Item(24, 1).Value = 7
$ws.Cells.Item(24, 2).Value = 180
$ws.Cells.Item(24, 3).Value = 'Casa de Cinema da Bahia'
$ws.Cells.Item(24, 4).Value = 'movie_theater'
$ws.Cells.Item(24, 5).Value = 1.248761228523516
$ws.Cells.Item(24, 6).Value = 0.642
$ws.Cells.Item(24, 7).Value = -12.9725511
$ws.Cells.Item(24, 8).Value = -38.508565

$ws.Cells.Item(25, 1).Value = 8
$ws.Cells.Item(25, 2).Value = 352
$ws.Cells.Item(25, 3).Value = 'Cia Maritima Salvador Shopping'
$ws.Cells.Item(25, 4).Value = 'shopping_mall'
$ws.Cells.Item(25, 5).Value = 1.279649308173254
$ws.Cells.Item(25, 6).Value = 0.677
$ws.Cells.Item(25, 7).Value = -12.9954162
$ws.Cells.Item(25, 8).Value = -38.4542581

$ws.Cells.Item(26, 1).Value = 9
$ws.Cells.Item(26, 2).Value = 83
$ws.Cells.Item(26, 3).Value = 'Ondina Beach'
$ws.Cells.Item(26, 4).Value = 'beach'
$ws.Cells.Item(26, 5).Value = 1.245057115085724
$ws.Cells.Item(26, 6).Value = 0.66
$ws.Cells.Item(26, 7).Value = -13.0103784
$ws.Cells.Item(26, 8).Value = -38.5120368

$ws.Cells.Item(27, 1).Value = 'NDCG: 0.9960169407355046'
$ws.Cells.Item(29, 1).Value = 'Diversificado_recs_random AWM'

# Table 3 (rows 31-40): Diversificado_recs_random AWM results
$ws.Cells.Item(31, 1).Value = 0
$ws.Cells.Item(31, 2).Value = 158
$ws.Cells.Item(31, 3).Value = 'Cazolla Gastro Burguer Beer.'
$ws.Cells.Item(31, 4).Value = 'fast-food'
$ws.Cells.Item(31, 5).Value = 1
$ws.Cells.Item(31, 6).Value = 0.516
$ws.Cells.Item(31, 7).Value = -12.9896638
$ws.Cells.Item(31, 8).Value = -38.4596392

$ws.Cells.Item(32, 1).Value = 1
$ws.Cells.Item(32, 2).Value = 50
$ws.Cells.Item(32, 3).Value = 'Bar Lagoa dos Frades'
$ws.Cells.Item(32, 4).Value = 'bar'
$ws.Cells.Item(32, 5).Value = 1
$ws.Cells.Item(32, 6).Value = 0.526
$ws.Cells.Item(32, 7).Value = -12.9807799
$ws.Cells.Item(32, 8).Value = -38.4436399

$ws.Cells.Item(33, 1).Value = 2
$ws.Cells.Item(33, 2).Value = 87
$ws.Cells.Item(33, 3).Value = 'Praia do Farol'
$ws.Cells.Item(33, 4).Value = 'beach'
$ws.Cells.Item(33, 5).Value = 1.581773681253018
$ws.Cells.Item(33, 6).Value = 0.791
$ws.Cells.Item(33, 7).Value = -12.952961
$ws.Cells.Item(33, 8).Value = -38.348126

$ws.Cells.Item(34, 1).Value = 3
$ws.Cells.Item(34, 2).Value = 266
$ws.Cells.Item(34, 3).Value = 'Beach Stop'
$ws.Cells.Item(34, 4).Value = 'restaurant'
$ws.Cells.Item(34, 5).Value = 1
$ws.Cells.Item(34, 6).Value = 0.528
$ws.Cells.Item(34, 7).Value = -12.932848
$ws.Cells.Item(34, 8).Value = -38.3287828

$ws.Cells.Item(35, 1).Value = 4
$ws.Cells.Item(35, 2).Value = 59
$ws.Cells.Item(35, 3).Value = 'Nova Alegria Bar e Restaurante'
$ws.Cells.Item(35, 4).Value = 'bar'
$ws.Cells.Item(35, 5).Value = 1
$ws.Cells.Item(35, 6).Value = 0.517
$ws.Cells.Item(35, 7).Value = -12.9693386
$ws.Cells.Item(35, 8).Value = -38.437136

$ws.Cells.Item(36, 1).Value = 5
$ws.Cells.Item(36, 2).Value = 377
$ws.Cells.Item(36, 3).Value = 'Letreiro Salvador'
$ws.Cells.Item(36, 4).Value = 'tourist_attraction'
$ws.Cells.Item(36, 5).Value = 0.1221469833687715
$ws.Cells.Item(36, 6).Value = 0.08699999999999999
$ws.Cells.Item(36, 7).Value = -12.974585
$ws.Cells.Item(36, 8).Value = -38.5127509

$ws.Cells.Item(37, 1).Value = 6
$ws.Cells.Item(37, 2).Value = 140
$ws.Cells.Item(37, 3).Value = 'CAFETERIA FLORIDA'
$ws.Cells.Item(37, 4).Value = 'cafe'
$ws.Cells.Item(37, 5).Value = 1
$ws.Cells.Item(37, 6).Value = 0.5
$ws.Cells.Item(37, 7).Value = -13.0039875
$ws.Cells.Item(37, 8).Value = -38.5023199

$ws.Cells.Item(38, 1).Value = 7
$ws.Cells.Item(38, 2).Value = 96
$ws.Cells.Item(38, 3).Value = 'Praia de Stella Maris'
$ws.Cells.Item(38, 4).Value = 'beach'
$ws.Cells.Item(38, 5).Value = 1.145291351242065
$ws.Cells.Item(38, 6).Value = 0.573
$ws.Cells.Item(38, 7).Value = -12.9488167
$ws.Cells.Item(38, 8).Value = -38.341097

$ws.Cells.Item(39, 1).Value = 8
$ws.Cells.Item(39, 2).Value = 60
$ws.Cells.Item(39, 3).Value = 'BaO Petiscaria'
$ws.Cells.Item(39, 4).Value = 'bar'
$ws.Cells.Item(39, 5).Value = 1
$ws.Cells.Item(39, 6).Value = 0.517
$ws.Cells.Item(39, 7).Value = -12.97853
$ws.Cells.Item(39, 8).Value = -38.4450989

$ws.Cells.Item(40, 1).Value = 9
$ws.Cells.Item(40, 2).Value = 93
$ws.Cells.Item(40, 3).Value = 'Praia da Boa Viagem'
$ws.Cells.Item(40, 4).Value = 'beach'
$ws.Cells.Item(40, 5).Value = 1.129106459124083
$ws.Cells.Item(40, 6).Value = 0.5649999999999999
$ws.Cells.Item(40, 7).Value = -12.9350958
$ws.Cells.Item(40, 8).Value = -38.5098867

$ws.Cells.Item(41, 1).Value = 'NDCG: 0.9247516059377486'
